# Presupuesto.xlsx revision edit
# - Strip the stray leading/trailing spaces that had crept into the recurring
#   " & " / " \\ \hline" / " \\ \cline{2-5}" LaTeX-row-separator text used across
#   the budget rows' B/D/F/I/L helper columns.
# - Bump the unit price of "Computador Portátil" (row 9) from 800 to 1000.
# - Extend the "Material de oficina" subtotal (N14) to include row 9, since it
#   now belongs in that block's sum.
# - Move the active selection to G25 (where the reviewer was last looking).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that carry the recurring CONCATENATE helper cells (B/D/F/I use "&",
# L uses either "\\ \hline" or "\\ \cline{2-5}").
$rowsHline = @(2,7,13,16,19,23,25)
$rowsCline = @(4,5,6,9,10,11,12,15,18,21,22)

foreach ($r in $rowsHline) {
    $ws.Cells.Item($r, 2).Value = "&"
    $ws.Cells.Item($r, 4).Value = "&"
    $ws.Cells.Item($r, 6).Value = "&"
    $ws.Cells.Item($r, 9).Value = "&"
    $ws.Cells.Item($r, 12).Value = "\\ \hline"
}

foreach ($r in $rowsCline) {
    $ws.Cells.Item($r, 2).Value = "&"
    $ws.Cells.Item($r, 4).Value = "&"
    $ws.Cells.Item($r, 6).Value = "&"
    $ws.Cells.Item($r, 9).Value = "&"
    $ws.Cells.Item($r, 12).Value = "\\ \cline{2-5}"
}

# Computador Portátil: unit price 800 -> 1000
$ws.Range("G9").Value = 1000

# The "Material de oficina" subtotal now spans rows 9-13 (was 10-13)
$ws.Range("N14").Formula = "=+SUM(J9:J13)"

# Reviewer's last-selected cell
$ws.Range("G25").Select() | Out-Null
